# Sample Project / Main.xlsx — "Rules" sheet
# Cell B11 ("R40") is replaced with the text value "1" (kept as text, not
# converted to a number — the cell must remain a shared-string cell so the
# existing General number format / style on B11 is left untouched).
#
# Setting .Value = "1" directly would make Excel auto-coerce the literal to
# a genuine number (losing the text type), so instead we briefly drop in a
# text-producing formula and immediately flatten it back down to a plain
# value with Copy / Paste-Special-Values, exactly like a user would do via
# Paste Special ▸ Values after typing a formula. This preserves B11's
# original cell style/number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("B11")
$target.Formula = '="1"'
$target.Copy()
$target.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
